# Update "Ready for handoff" status text to "In Translation" wherever it
# occurs (Overview!E2:F2, zh-cn!C2, de-de!C2 all share this string).
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: keep the literal on the left of -eq; PowerShell coerces the
        # right-hand operand to the left-hand operand's type, and a bare
        # boolean cell value (e.g. TRUE/FALSE cells) would otherwise coerce
        # this string literal to $true and falsely match every truthy cell.
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value = "In Translation"
        }
    }
}

# Narrow the "Status" columns now that the shorter "In Translation" text no
# longer needs as much room (previously sized for "Ready for handoff"):
# Overview!E:F (zh-cn / de-de status columns) and the "Status" table column
# on each locale sheet (zh-cn!C, de-de!C).
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
